$d = $word.ActiveDocument

$d.Content.Find.Execute("615÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "498÷5=", 2)
$d.Content.Find.Execute("226÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "649÷4=", 2)
$d.Content.Find.Execute("976÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "906÷8=", 2)
$d.Content.Find.Execute("720÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "235÷2=", 2)
$d.Content.Find.Execute("657÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "591÷9=", 2)
$d.Content.Find.Execute("432÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "172÷5=", 2)
$d.Content.Find.Execute("347÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "990÷2=", 2)
$d.Content.Find.Execute("762÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "388÷6=", 2)
$d.Content.Find.Execute("219÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "802÷8=", 2)
$d.Content.Find.Execute("531÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "614÷3=", 2)
$d.Content.Find.Execute("841÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "924÷4=", 2)
$d.Content.Find.Execute("612÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "420÷5=", 2)
$d.Content.Find.Execute("876÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "699÷3=", 2)
$d.Content.Find.Execute("835÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "156÷7=", 2)
$d.Content.Find.Execute("986÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "386÷7=", 2)
$d.Content.Find.Execute("951÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "569÷5=", 2)
$d.Content.Find.Execute("626÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "371÷6=", 2)
$d.Content.Find.Execute("364÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "358÷2=", 2)
$d.Content.Find.Execute("319÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "909÷7=", 2)
$d.Content.Find.Execute("170÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "844÷9=", 2)
$d.Content.Find.Execute("119÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "605÷7=", 2)
$d.Content.Find.Execute("845÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "850÷4=", 2)
$d.Content.Find.Execute("220÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "439÷3=", 2)
$d.Content.Find.Execute("488÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "895÷2=", 2)
$d.Content.Find.Execute("760÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "124÷4=", 2)
